$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 22: image-name headers, merged in pairs (D:E, F:G, H:I, J:K)
# ---------------------------------------------------------------------------
$ws.Range("A22:C22").Style = "Total"

$ws.Range("D22").Value = "Horses_Run_Animals_horse_9192x6012"
$ws.Range("F22").Value = "Red_Mazda_2528_1368"
$ws.Range("H22").Value = "range_rover_1920_1080"
$ws.Range("J22").Value = "daimler_800_777"

$ws.Range("D22:K22").Style = "Total"
$ws.Range("D22:K22").HorizontalAlignment = -4108

$ws.Range("D22:E22").Merge()
$ws.Range("F22:G22").Merge()
$ws.Range("H22:I22").Merge()
$ws.Range("J22:K22").Merge()

$ws.Rows.Item(22).RowHeight = 15.75

# ---------------------------------------------------------------------------
# Row 23: CPU / GPU sub-headers
# ---------------------------------------------------------------------------
$ws.Range("A23:C23").Style = "Total"

$ws.Range("D23").Value = "CPU"
$ws.Range("E23").Value = "GPU"
$ws.Range("F23").Value = "CPU"
$ws.Range("G23").Value = "GPU"
$ws.Range("H23").Value = "CPU"
$ws.Range("I23").Value = "GPU"
$ws.Range("J23").Value = "CPU"
$ws.Range("K23").Value = "GPU"

$ws.Range("D23:K23").Style = "Total"

$ws.Rows.Item(23).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Rows 24-32: data
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "apply_gaussian_kernel"
$ws.Range("A25").Value = "compute_pixel_thresholds"
$ws.Range("A26").Value = "apply_sobel_filter_x"
$ws.Range("A27").Value = "apply_sobel_filter_y"
$ws.Range("A28").Value = "calculate_sobel_magnitude"
$ws.Range("A29").Value = "calculate_sobel_direction"
$ws.Range("A30").Value = "apply_non_max_suppression"
$ws.Range("A31").Value = "apply_double_thresholds"
$ws.Range("A32").Value = "apply_hysteresis_edge_tracking"

$ws.Rows.Item(24).RowHeight = 15.75

$ws.Range("D24").Value = 1629.69
$ws.Range("E24").Value = 17.352029999999999
$ws.Range("F24").Value = 103.37
$ws.Range("G24").Value = 1.07213
$ws.Range("H24").Value = 62.06
$ws.Range("I24").Value = 0.64866999999999997
$ws.Range("J24").Value = 18.739999999999998
$ws.Range("K24").Value = 0.20058000000000001

$ws.Range("D25").Value = 63.61
$ws.Range("E25").Value = 12.748799999999999
$ws.Range("F25").Value = 3.98
$ws.Range("G25").Value = 0.87002000000000002
$ws.Range("H25").Value = 2.39
$ws.Range("I25").Value = 0.51500999999999997
$ws.Range("J25").Value = 0.72
$ws.Range("K25").Value = 0.21187

$ws.Range("D26").Value = 3192.68
$ws.Range("E26").Value = 0.00163
$ws.Range("F26").Value = 197.6
$ws.Range("G26").Value = 0.0017
$ws.Range("H26").Value = 122.05
$ws.Range("I26").Value = 0.0017
$ws.Range("J26").Value = 35.380000000000003
$ws.Range("K26").Value = 0.0008

$ws.Range("D27").Value = 3135.95
$ws.Range("E27").Value = 0.06758
$ws.Range("F27").Value = 194.29
$ws.Range("G27").Value = 0.04304
$ws.Range("H27").Value = 127.96
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 34.74
$ws.Range("K27").Value = 0

$ws.Range("D28").Value = 254.7
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 16.149999999999999
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 11.81
$ws.Range("I28").Value = 0.03059
$ws.Range("J28").Value = 3.01
$ws.Range("K28").Value = 0.03542

$ws.Range("D29").Value = 561.04
$ws.Range("E29").Value = 5.3512599999999999
$ws.Range("F29").Value = 33.979999999999997
$ws.Range("G29").Value = 0.33821000000000001
$ws.Range("H29").Value = 20.36
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 6.32
$ws.Range("K29").Value = 0

$ws.Range("D30").Value = 420.45
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 21.3
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 15.89
$ws.Range("I30").Value = 0.2089
$ws.Range("J30").Value = 4.29
$ws.Range("K30").Value = 0.06253

$ws.Range("D31").Value = 195.8
$ws.Range("E31").Value = 2.8958699999999999
$ws.Range("F31").Value = 10.35
$ws.Range("G31").Value = 0.18890000000000001
$ws.Range("H31").Value = 6
$ws.Range("I31").Value = 0.11165
$ws.Range("J31").Value = 1.92
$ws.Range("K31").Value = 0.03565

$ws.Range("D32").Value = 209.98
$ws.Range("E32").Value = 8.4365400000000008
$ws.Range("F32").Value = 10.3
$ws.Range("G32").Value = 0.42249999999999999
$ws.Range("H32").Value = 6.4
$ws.Range("I32").Value = 0.25600000000000001
$ws.Range("J32").Value = 1.83
$ws.Range("K32").Value = 0.07024

# ---------------------------------------------------------------------------
# Row 33: Total
# ---------------------------------------------------------------------------
$ws.Range("A33").Value = "Total"
$ws.Range("A33:K33").Style = "Calculation"

$ws.Range("D33").Formula = "=SUM(D24:D32)"
$ws.Range("E33").Formula = "=SUM(E24:E32)"
$ws.Range("F33").Formula = "=SUM(F24:F32)"
$ws.Range("G33").Formula = "=SUM(G24:G32)"
$ws.Range("H33").Formula = "=SUM(H24:H32)"
$ws.Range("I33").Formula = "=SUM(I24:I32)"
$ws.Range("J33").Formula = "=SUM(J24:J32)"
$ws.Range("K33").Formula = "=SUM(K24:K32)"

# ---------------------------------------------------------------------------
# Row 34: Speed V/S CPU
# ---------------------------------------------------------------------------
$ws.Range("A34").Value = "Speed V/S CPU"
$ws.Range("A34:K34").Style = "Calculation"

$ws.Range("D34").Formula = "=D33/D33"
$ws.Range("E34").Formula = "=D33/E33"
$ws.Range("F34").Formula = "=F33/F33"
$ws.Range("G34").Formula = "=F33/G33"
$ws.Range("H34").Formula = "=H33/H33"
$ws.Range("I34").Formula = "=H33/I33"
$ws.Range("J34").Formula = "=J33/J33"
$ws.Range("K34").Formula = "=J33/K33"

# ---------------------------------------------------------------------------
# Column widths for the new columns
# ---------------------------------------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 14.8
$ws.Columns.Item(8).ColumnWidth = 14.1
$ws.Columns.Item(9).ColumnWidth = 12
$ws.Columns.Item(10).ColumnWidth = 9.6
$ws.Columns.Item(11).ColumnWidth = 13.8

# ---------------------------------------------------------------------------
# Selection, to mirror the saved workbook view
# ---------------------------------------------------------------------------
$ws.Range("C17").Select()
